# Updates the Price (D) and Volume(1h) (E) columns of the cryptos list
# to the latest scraped values, as produced by the GitHub Actions job.
#
# Note: some new Price values (e.g. "1.00", "7.75") would otherwise be
# auto-converted by Excel into numbers, losing the trailing-zero text
# formatting used throughout this sheet. For those cells we assign the
# value through .Formula using a leading single-quote (the same trick
# Excel itself uses to force a value to be stored as literal text),
# which keeps the cell's text representation exactly as scraped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.284.00'
$ws.Range('E2').Value = '  +1.83%  '
# Row 3
$ws.Range('D3').Value = '3.466.32'
$ws.Range('E3').Value = '  +1.13%  '
# Row 4
$ws.Range('D4').Formula = '''1.00'
$ws.Range('E4').Value = '  +0.03%  '
# Row 5
$ws.Range('D5').Formula = '''581.50'
$ws.Range('E5').Value = '  +0.42%  '
# Row 6
$ws.Range('D6').Formula = '''147.54'
$ws.Range('E6').Value = '  +1.66%  '
# Row 7
$ws.Range('D7').Value = '3.465.75'
$ws.Range('E7').Value = '  +1.13%  '
# Row 8
$ws.Range('E8').Value = '  -0.05%  '
# Row 9
$ws.Range('E9').Value = '  +0.43%  '
# Row 10
$ws.Range('D10').Formula = '''7.75'
$ws.Range('E10').Value = '  +1.70%  '
# Row 11
$ws.Range('E11').Value = '  +0.68%  '
# Row 12
$ws.Range('E12').Value = '  +4.24%  '
# Row 13
$ws.Range('D13').Value = '4.060.53'
$ws.Range('E13').Value = '  +1.13%  '
# Row 14
$ws.Range('D14').Formula = '''29.50'
$ws.Range('E14').Value = '  +2.39%  '
# Row 15
$ws.Range('E15').Value = '  +2.58%  '
# Row 16
$ws.Range('D16').Value = '3.472.85'
$ws.Range('E16').Value = '  +1.19%  '
# Row 17
$ws.Range('E17').Value = '  +0.59%  '
# Row 18
$ws.Range('D18').Value = '63.294.71'
$ws.Range('E18').Value = '  +1.72%  '
# Row 19
$ws.Range('E19').Value = '  +3.04%  '
# Row 20
$ws.Range('D20').Formula = '''14.51'
$ws.Range('E20').Value = '  +3.26%  '
# Row 21
$ws.Range('E21').Value = '  +1.36%  '
# Row 22
$ws.Range('D22').Formula = '''389.03'
$ws.Range('E22').Value = '  -1.24%  '
# Row 23
$ws.Range('E23').Value = '  +1.75%  '
# Row 24
$ws.Range('D24').Formula = '''74.47'
$ws.Range('E24').Value = '  -0.58%  '
# Row 25
$ws.Range('E25').Value = '  -0.17%  '
# Row 26
$ws.Range('D26').Value = '3.614.27'
$ws.Range('E26').Value = '  +1.29%  '
# Row 27
$ws.Range('E27').Value = '  +0.04%  '
# Row 28
$ws.Range('E28').Value = '  -2.89%  '
# Row 29
$ws.Range('E29').Value = '  +1.65%  '
# Row 30
$ws.Range('D30').Formula = '''0.999'
$ws.Range('E30').Value = '  +0.15%  '
# Row 31
$ws.Range('D31').Formula = '''8.21'
# Row 32
$ws.Range('E32').Value = '  -0.31%  '
# Row 33
$ws.Range('E33').Value = '  -0.01%  '
# Row 34
$ws.Range('E34').Value = '  -4.27%  '
# Row 35
$ws.Range('E35').Value = '  -0.78%  '
# Row 36
$ws.Range('D36').Formula = '''1.63'
$ws.Range('E36').Value = '  +5.30%  '
# Row 37
$ws.Range('E37').Value = '  +0.01%  '
# Row 38
$ws.Range('D38').Formula = '''7.15'
$ws.Range('E38').Value = '  +1.93%  '
# Row 39
$ws.Range('D39').Formula = '''31.93'
$ws.Range('E39').Value = '  +10.49%  '
# Row 40
$ws.Range('D40').Formula = '''167.85'
$ws.Range('E40').Value = '  +0.05%  '
# Row 41
$ws.Range('D41').Value = '3.504.60'
$ws.Range('E41').Value = '  +1.25%  '
# Row 42
$ws.Range('E42').Value = '  +1.36%  '
# Row 43
$ws.Range('D43').Formula = '''0.793'
$ws.Range('E43').Value = '  +0.40%  '
# Row 44
$ws.Range('E44').Value = '  +3.72%  '
# Row 45
$ws.Range('E45').Value = '  -1.12%  '
# Row 46
$ws.Range('E46').Value = '  +3.05%  '
# Row 47
$ws.Range('E47').Value = '  -1.52%  '
# Row 48
$ws.Range('D48').Value = '2.593.69'
$ws.Range('E48').Value = '  +3.09%  '
# Row 49
$ws.Range('D49').Formula = '''2.33'
$ws.Range('E49').Value = '  +8.41%  '
# Row 50
$ws.Range('E50').Value = '  +1.99%  '
# Row 51
$ws.Range('D51').Formula = '''23.02'
$ws.Range('E51').Value = '  -0.42%  '
